$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)

# PowerPoint COM measures shape position/size in points; the OOXML stores
# EMU (914400 EMU/inch, 12700 EMU/point). Add a tiny epsilon before
# converting so that float round-trip error doesn't truncate the EMU value
# down by one.
$EMU_PER_PT = 12700
$epsilon = 0.25 / $EMU_PER_PT

$shape.Left = (178944 / $EMU_PER_PT) + $epsilon
$shape.Top = (922010 / $EMU_PER_PT) + $epsilon
$shape.Width = (3406877 / $EMU_PER_PT) + $epsilon
$shape.Height = (400110 / $EMU_PER_PT) + $epsilon

# Replace the paragraph's two runs ("{{titulo}}{{resumo}} {{data}}{{link}}"
# followed by a lone space run) with a single run containing just
# "{{titulo}}".
$shape.TextFrame.TextRange.Text = "{{titulo}}"
